$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data
$ws.Range("D2").Value = "67.934.57"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "3.329.23"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").Value = "3.325.89"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("E10").Value = "  +5.59%  "
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.68%  "
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "683.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "3.870.83"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.91%  "
$ws.Range("D17").Value = "67.942.42"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").Value = "3.330.76"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.896"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("E30").Value = "  +2.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "565.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.00%  "
$ws.Range("E34").Value = "  +3.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "3.704.08"
$ws.Range("E37").Value = "  -4.31%  "
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.17%  "
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("E41").Value = "  +6.82%  "
$ws.Range("E42").Value = "  +2.47%  "
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").Value = "0.0₃0674"
$ws.Range("E44").Value = "  +1.19%  "
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0406"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("E47").Value = "  +5.40%  "
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  -3.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.36%  "
